$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.665.31"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.29"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.63"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4465"
$ws.Range("E7").Value = "  +4.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("E8").Value = "  +8.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.90"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.157"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07620"
$ws.Range("E11").Value = "  +2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.89"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9983"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.377"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.627"
$ws.Range("E15").Value = "  +4.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.801.61"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001099"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06756"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.54"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9982"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.88"
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.364"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.666.63"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.414"
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.80"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.394"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.52"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.006.97"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.20"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.277"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.956"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.909"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09382"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2297"
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.39"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06415"
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02353"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.247"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6668"
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.240"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.456"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9979"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6172"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.823"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.83"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.073"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07118"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.164"
$ws.Range("E51").Value = "  -1.72%  "
